# CI: Update Excel counters (state_counters + packages)
#
# Appends 13 new "Packages" rows (rows 9-21) to the Packages worksheet,
# matching new CMS package/status reference data added in the source
# commit. Columns are:
#   A=PackageType  B=State  C=Authority  D=ActionType
#   E=PackageID    F=Status G=ParentID

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9 : SPA / MD / Medicaid SPA / MD-25-9473 / Under Review ---
$ws.Range("A9").Value = 'SPA'
$ws.Range("B9").Value = 'MD'
$ws.Range("C9").Value = 'Medicaid SPA'
$ws.Range("E9").Value = 'MD-25-9473'
$ws.Range("F9").Value = 'Under Review'

# --- Row 10 : SPA / MD / Medicaid SPA / MD-25-9474 / Approved ---
$ws.Range("A10").Value = 'SPA'
$ws.Range("B10").Value = 'MD'
$ws.Range("C10").Value = 'Medicaid SPA'
$ws.Range("E10").Value = 'MD-25-9474'
$ws.Range("F10").Value = 'Approved'

# --- Row 11 : SPA / MD / Medicaid SPA / MD-25-9475 / Submitted ---
$ws.Range("A11").Value = 'SPA'
$ws.Range("B11").Value = 'MD'
$ws.Range("C11").Value = 'Medicaid SPA'
$ws.Range("E11").Value = 'MD-25-9475'
$ws.Range("F11").Value = 'Submitted'

# --- Row 12 : SPA / MD / Medicaid SPA / MD-25-9476 / Disapproved ---
$ws.Range("A12").Value = 'SPA'
$ws.Range("B12").Value = 'MD'
$ws.Range("C12").Value = 'Medicaid SPA'
$ws.Range("E12").Value = 'MD-25-9476'
$ws.Range("F12").Value = 'Disapproved'

# --- Row 13 : Waiver / MD / 1915(c) / Amendment / MD-2260.R00.42 / (no status) / MD-2260.R00.00 ---
$ws.Range("A13").Value = 'Waiver'
$ws.Range("B13").Value = 'MD'
$ws.Range("C13").Value = '1915(c)'
$ws.Range("D13").Value = 'Amendment'
$ws.Range("E13").Value = 'MD-2260.R00.42'
$ws.Range("G13").Value = 'MD-2260.R00.00'

# --- Row 14 : SPA / MD / Medicaid SPA / MD-25-9477 / Pending-Concurrence ---
$ws.Range("A14").Value = 'SPA'
$ws.Range("B14").Value = 'MD'
$ws.Range("C14").Value = 'Medicaid SPA'
$ws.Range("E14").Value = 'MD-25-9477'
$ws.Range("F14").Value = 'Pending-Concurrence'

# --- Row 15 : SPA / MD / Medicaid SPA / MD-25-9478 / RAI Issued ---
$ws.Range("A15").Value = 'SPA'
$ws.Range("B15").Value = 'MD'
$ws.Range("C15").Value = 'Medicaid SPA'
$ws.Range("E15").Value = 'MD-25-9478'
$ws.Range("F15").Value = 'RAI Issued'

# --- Row 16 : SPA / MD / CHIP SPA / MD-25-9479 / Submitted ---
$ws.Range("A16").Value = 'SPA'
$ws.Range("B16").Value = 'MD'
$ws.Range("C16").Value = 'CHIP SPA'
$ws.Range("E16").Value = 'MD-25-9479'
$ws.Range("F16").Value = 'Submitted'

# --- Row 17 : SPA / MD / Medicaid SPA / MD-25-9480 / Submitted ---
$ws.Range("A17").Value = 'SPA'
$ws.Range("B17").Value = 'MD'
$ws.Range("C17").Value = 'Medicaid SPA'
$ws.Range("E17").Value = 'MD-25-9480'
$ws.Range("F17").Value = 'Submitted'

# --- Row 18 : Waiver / MD / 1915(c) / Amendment / MD-2260.R00.43 / Unsubmitted / MD-2260.R00.00 ---
$ws.Range("A18").Value = 'Waiver'
$ws.Range("B18").Value = 'MD'
$ws.Range("C18").Value = '1915(c)'
$ws.Range("D18").Value = 'Amendment'
$ws.Range("E18").Value = 'MD-2260.R00.43'
$ws.Range("F18").Value = 'Unsubmitted'
$ws.Range("G18").Value = 'MD-2260.R00.00'

# --- Row 19 : SPA / MD / Medicaid SPA / MD-25-9481 / Under Review ---
$ws.Range("A19").Value = 'SPA'
$ws.Range("B19").Value = 'MD'
$ws.Range("C19").Value = 'Medicaid SPA'
$ws.Range("E19").Value = 'MD-25-9481'
$ws.Range("F19").Value = 'Under Review'

# --- Row 20 : Waiver / MD / 1915(b) / Initial / MD-2278.R00.00 / Terminated ---
$ws.Range("A20").Value = 'Waiver'
$ws.Range("B20").Value = 'MD'
$ws.Range("C20").Value = '1915(b)'
$ws.Range("D20").Value = 'Initial'
$ws.Range("E20").Value = 'MD-2278.R00.00'
$ws.Range("F20").Value = 'Terminated'

# --- Row 21 : SPA / MD / Medicaid SPA / MD-25-9482 / Withdrawn ---
$ws.Range("A21").Value = 'SPA'
$ws.Range("B21").Value = 'MD'
$ws.Range("C21").Value = 'Medicaid SPA'
$ws.Range("E21").Value = 'MD-25-9482'
$ws.Range("F21").Value = 'Withdrawn'

# Cells that are blank in the source data (ActionType/ParentID for SPA rows,
# Status for the bare waiver-amendment row) still need to be real empty-text
# cells (matching the existing empty-string shared-string entry used
# elsewhere in this sheet), not entirely-unset cells. A plain
# `Range.Value = ""` leaves the cell completely empty (no cell record at
# all), so instead we enter a lone `'` (Excel's text quote-prefix), which
# forces an empty-text cell, then reset the cell style back to Normal so we
# don't leave a stray quote-prefixed style behind.
$emptyCells = @(
    "D9","G9",
    "D10","G10",
    "D11","G11",
    "D12","G12",
    "F13",
    "D14","G14",
    "D15","G15",
    "D16","G16",
    "D17","G17",
    "D19","G19",
    "G20",
    "D21","G21"
)
foreach ($addr in $emptyCells) {
    $ws.Range($addr).Value = "'"
}
foreach ($addr in $emptyCells) {
    $ws.Range($addr).Style = "Normal"
}
